# Update the FeSources sheet: turn the plain URL cells in column D into
# HTML anchor-tag strings (<a href='...'>Label</a>), and fix the
# "Working Futures 2035" label to "Working Futures 20235" in A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = "<a href='https://www.ons.gov.uk/peoplepopulationandcommunity/wellbeing/articles/subnationalindicatorsexplorer/2022-01-06'>ONS</a>"
$ws.Range("D4").Value = "<a href='https://www.aoc.co.uk/research-unit/data-sources'>AOC</a>"
$ws.Range("D5").Value = "<a href='https://census.gov.uk/local-authorities'>Census</a>"
$ws.Range("D6").Value = "<a href='https://www.nomisweb.co.uk/'>Nomis</a>"
$ws.Range("D7").Value = "<a href='https://explore-education-statistics.service.gov.uk/'>EES</a>"

$ws.Range("A8").Value = "Working Futures 20235"

# Move/restore the sheet selection to E10, matching the saved view state.
$ws.Range("E10").Select()
